$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 174
$data = @(
    @("08-09-2021", 14575, 19545, -4970),
    @("09-09-2021", 14530, 19836, -5306),
    @("10-09-2021", 14577, 19847, -5270),
    @("13-09-2021", 14531, 19990, -5458),
    @("14-09-2021", 14731, 19814, -5084),
    @("15-09-2021", 14719, 19591, -4872),
    @("16-09-2021", 14891, 19655, -4764),
    @("20-09-2021", 14729, 19467, -4738),
    @("21-09-2021", 14441, 19384, -4943),
    @("22-09-2021", 14736, 19962, -5225),
    @("23-09-2021", 14423, 19995, -5572),
    @("24-09-2021", 14401, 19837, -5436),
    @("27-09-2021", 14384, 19950, -5566),
    @("28-09-2021", 13421, 19481, -6060),
    @("29-09-2021", 13503, 19680, -6177),
    @("30-09-2021", 13927, 19624, -5698)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $cellA = $ws.Cells.Item($row, 1)
    # The "Serie" column holds plain date-as-text labels (shared strings),
    # not real dates. A few of the new values (day-of-month <= 12) are
    # ambiguous and would otherwise get auto-recognized as a date value by
    # Excel's smart entry, so force those to text before writing them.
    $dayPart = [int]($rowData[0].Substring(0, 2))
    if ($dayPart -le 12) {
        $cellA.NumberFormat = "@"
        $cellA.Value = $rowData[0]
        $cellA.Style = "Normal"
    } else {
        $cellA.Value = $rowData[0]
    }

    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
